$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws.Range("B9").Value = "Alvearie Team"

# Replace the (duplicate) "Contact" / "No display for ContactDetail" row
# with a single "Jurisdiction" / "United States of America" row, and drop
# the now-redundant duplicate row entirely.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Rows.Item(11).Delete()
